$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell updates as described by the diff.
# Each cell is set to Text format first to guarantee the literal string
# is preserved (avoids Excel auto-converting numeric-looking strings like
# "88.45" or "0.999" into real numbers), then style is reset back to
# "Normal" so no stray style index is left on the cell.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "40.833.38"
Set-TextValue "E2" "  -2.05%  "
Set-TextValue "D3" "2.384.67"
Set-TextValue "E3" "  -3.57%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "313.15"
Set-TextValue "E5" "  -1.80%  "
Set-TextValue "D6" "88.45"
Set-TextValue "E6" "  -4.89%  "
Set-TextValue "D7" "0.531"
Set-TextValue "E7" "  -3.88%  "
Set-TextValue "E8" "  +0.01%  "
Set-TextValue "D9" "0.494"
Set-TextValue "E9" "  -4.53%  "
Set-TextValue "D10" "0.0825"
Set-TextValue "E10" "  -4.63%  "
Set-TextValue "D11" "31.25"
Set-TextValue "E11" "  -6.06%  "
Set-TextValue "E12" "  -1.73%  "
Set-TextValue "D13" "2.750.98"
Set-TextValue "E13" "  -3.73%  "
Set-TextValue "D14" "6.58"
Set-TextValue "E14" "  -4.56%  "
Set-TextValue "D15" "15.18"
Set-TextValue "E15" "  -3.51%  "
Set-TextValue "D16" "2.384.38"
Set-TextValue "E16" "  -4.54%  "
Set-TextValue "D17" "0.765"
Set-TextValue "E17" "  -3.49%  "
Set-TextValue "D18" "40.733.32"
Set-TextValue "E18" "  -2.20%  "
Set-TextValue "B19" "ShibaInu"
Set-TextValue "C19" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D19" "0.0₃0914"
Set-TextValue "E19" "  -3.83%  "
Set-TextValue "B20" "Uniswap"
Set-TextValue "C20" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "6.21"
Set-TextValue "E20" "  -3.96%  "
Set-TextValue "D21" "69.44"
Set-TextValue "E21" "  -2.40%  "
Set-TextValue "D22" "10.85"
Set-TextValue "E22" "  -3.89%  "
Set-TextValue "D23" "234.08"
Set-TextValue "E23" "  -2.17%  "
Set-TextValue "E24" "  -3.33%  "
Set-TextValue "E25" "  +0.11%  "
Set-TextValue "E26" "  -6.02%  "
Set-TextValue "D27" "23.79"
Set-TextValue "E27" "  -3.52%  "
Set-TextValue "E28" "  -2.29%  "
Set-TextValue "D29" "9.40"
Set-TextValue "E29" "  -4.02%  "
Set-TextValue "D30" "33.97"
Set-TextValue "E30" "  -5.78%  "
Set-TextValue "D31" "155.78"
Set-TextValue "E31" "  -1.97%  "
Set-TextValue "E32" "  -0.13%  "
Set-TextValue "E33" "  -5.41%  "
Set-TextValue "D34" "0.0733"
Set-TextValue "E34" "  -4.31%  "
Set-TextValue "D35" "2.43"
Set-TextValue "E35" "  -6.27%  "
Set-TextValue "D36" "0.114"
Set-TextValue "E36" "  -2.07%  "
Set-TextValue "B37" "LidoDAOToken"
Set-TextValue "C37" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D37" "2.82"
Set-TextValue "E37" "  -3.62%  "
Set-TextValue "B38" "Celestia"
Set-TextValue "C38" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D38" "16.16"
Set-TextValue "E38" "  -7.90%  "
Set-TextValue "E39" "  -3.18%  "
Set-TextValue "E40" "  -7.26%  "
Set-TextValue "D41" "3.82"
Set-TextValue "E41" "  -5.28%  "
Set-TextValue "E42" "  -7.82%  "
Set-TextValue "D43" "1.954.47"
Set-TextValue "E43" "  -1.94%  "
Set-TextValue "D44" "0.0272"
Set-TextValue "E44" "  -4.67%  "
Set-TextValue "D45" "17.54"
Set-TextValue "E45" "  -6.76%  "
Set-TextValue "E46" "  -6.38%  "
Set-TextValue "D47" "9.39"
Set-TextValue "E47" "  -0.80%  "
Set-TextValue "D48" "2.617.31"
Set-TextValue "E48" "  -3.60%  "
Set-TextValue "B49" "Aave"
Set-TextValue "C49" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D49" "93.91"
Set-TextValue "E49" "  -3.53%  "
Set-TextValue "B50" "BitcoinSV"
Set-TextValue "C50" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D50" "73.18"
Set-TextValue "E50" "  -0.91%  "
Set-TextValue "D51" "51.00"
Set-TextValue "E51" "  -3.35%  "
